# Apply crypto price/volume updates produced by the scheduled scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.053.42"
$ws.Range("E2").Value = "  +8.41%  "

$ws.Range("D3").Value = "3.513.78"
$ws.Range("E3").Value = "  +11.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.92"
$ws.Range("E5").Value = "  +14.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "548.92"
$ws.Range("E6").Value = "  +6.44%  "

$ws.Range("D7").Value = "3.505.94"
$ws.Range("E7").Value = "  +11.75%  "

$ws.Range("E8").Value = "  +3.63%  "

$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.630"
$ws.Range("E10").Value = "  +6.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  +18.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.80"
$ws.Range("E12").Value = "  +6.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  +9.90%  "

$ws.Range("E14").Value = "  +5.60%  "

$ws.Range("D15").Value = "4.072.55"
$ws.Range("E15").Value = "  +11.70%  "

$ws.Range("D16").Value = "3.509.32"
$ws.Range("E16").Value = "  +11.90%  "

$ws.Range("E17").Value = "  +5.61%  "

$ws.Range("D18").Value = "66.976.40"
$ws.Range("E18").Value = "  +8.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.14"
$ws.Range("E19").Value = "  +7.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.95"
$ws.Range("E20").Value = "  +11.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.992"
$ws.Range("E21").Value = "  +3.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.42"
$ws.Range("E22").Value = "  +19.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.90"
$ws.Range("E23").Value = "  +6.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.55"
$ws.Range("E24").Value = "  +6.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.16"
$ws.Range("E25").Value = "  +7.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.13"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("E27").Value = "  +13.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.92"
$ws.Range("E28").Value = "  +8.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.86"
$ws.Range("E29").Value = "  +11.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.12"
$ws.Range("E30").Value = "  +8.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "654.69"
$ws.Range("E31").Value = "  +3.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.64"
$ws.Range("E32").Value = "  +5.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.68"
$ws.Range("E33").Value = "  +5.06%  "

$ws.Range("E34").Value = "  +7.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.44"
$ws.Range("E35").Value = "  +6.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.46"
$ws.Range("E36").Value = "  +6.70%  "

$ws.Range("D37").Value = "0.0₃0814"
$ws.Range("E37").Value = "  +20.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.388"
$ws.Range("E39").Value = "  +6.07%  "

$ws.Range("E40").Value = "  +14.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.28"
$ws.Range("E41").Value = "  +15.11%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").Value = "2.991.23"
$ws.Range("E43").Value = "  +4.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  +5.71%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.86"
$ws.Range("E45").Value = "  +14.97%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.36"
$ws.Range("E46").Value = "  +14.60%  "

$ws.Range("E47").Value = "  +8.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.70"
$ws.Range("E48").Value = "  +4.17%  "

$ws.Range("E49").Value = "  +7.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.69"
$ws.Range("E50").Value = "  +17.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.01"
$ws.Range("E51").Value = "  +6.52%  "
